# Validation data for model "2-1" (simple FNN, residue oversampling)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# --- Row 7 ---
$ws.Cells.Item(7, 2).Value = "2-1"
$ws.Cells.Item(7, 3).Value = "simple FNN, residue oversampling"
$ws.Cells.Item(7, 4).Value = 0
$ws.Cells.Item(7, 5).Value = 0.05
$ws.Cells.Item(7, 6).Value = 53886
$ws.Cells.Item(7, 7).Value = 5881
$ws.Cells.Item(7, 8).Value = 39922
$ws.Cells.Item(7, 9).Value = 641
$ws.Cells.Item(7, 10).Formula = "=F7/(F7+G7)"
$ws.Cells.Item(7, 11).Formula = "=F7/(F7+I7)"
$ws.Cells.Item(7, 12).Formula = "=2*((J7*K7)/(J7+K7))"
$ws.Cells.Item(7, 13).Formula = "=(H7*F7-G7*I7)/SQRT((H7+I7)*(G7+F7)*(H7+G7)*(I7+F7))"
$ws.Cells.Item(7, 14).Formula = "=(K7+(H7/(H7+G7)))/2"

# --- Row 8 ---
$ws.Cells.Item(8, 2).Value = "2-1"
$ws.Cells.Item(8, 3).Value = "simple FNN, residue oversampling"
$ws.Cells.Item(8, 4).Value = 4
$ws.Cells.Item(8, 5).Value = 0.75
$ws.Cells.Item(8, 6).Value = 57604
$ws.Cells.Item(8, 7).Value = 7009
$ws.Cells.Item(8, 8).Value = 40492
$ws.Cells.Item(8, 9).Value = 60
$ws.Cells.Item(8, 10).Formula = "=F8/(F8+G8)"
$ws.Cells.Item(8, 11).Formula = "=F8/(F8+I8)"
$ws.Cells.Item(8, 12).Formula = "=2*((J8*K8)/(J8+K8))"
$ws.Cells.Item(8, 13).Formula = "=(H8*F8-G8*I8)/SQRT((H8+I8)*(G8+F8)*(H8+G8)*(I8+F8))"
$ws.Cells.Item(8, 14).Formula = "=(K8+(H8/(H8+G8)))/2"

$ws.Range("B2").Select()
